$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SCORE (D2) and CLUSTER MUST HAVE SCORE (G2) are stored as text in the
# workbook (e.g. "32.27"), so force text formatting before assigning the
# new values to stop Excel from auto-coercing them into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.13"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "31.2"

# CLUSTER MUST HAVE MATCH column (F2:F17) was re-ordered.
$ws.Range("F2").Value = "data architecture : 1"
$ws.Range("F3").Value = "analysis : 6"
$ws.Range("F4").Value = "access : 1"
$ws.Range("F5").Value = "engineer : 1"
$ws.Range("F6").Value = "database : 1"
$ws.Range("F7").Value = "hadoop : 1"
$ws.Range("F8").Value = "design : 4"
$ws.Range("F9").Value = "big data : 3"
$ws.Range("F10").Value = "aws : 2"
$ws.Range("F11").Value = "documentation : 2"
$ws.Range("F12").Value = "sql : 4"
$ws.Range("F13").Value = "tools : 2"
$ws.Range("F14").Value = "spark : 1"
$ws.Range("F15").Value = "analyze : 2"
$ws.Range("F16").Value = "python : 1"
$ws.Range("F17").Value = "data engineer : 1"

# CLUSTER GOOD TO HAVE MATCH column (H2:H3) values swapped.
$ws.Range("H2").Value = "aws : 2"
$ws.Range("H3").Value = "management : 2"
